$d = $word.ActiveDocument

function New-FlatOpc($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1. Title paragraph: "Votre proposition commerciale" -> split runs with
#    spell-check proofErr markers.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(8)
$titleBody = '<w:body><w:p><w:pPr><w:pStyle w:val="Titre"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Votre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> proposition </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>commerciale</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body>'
$titlePara.Range.InsertXML((New-FlatOpc $titleBody))

# ---------------------------------------------------------------------------
# 2. Replace the big block of paragraphs 10..25 (after renumbering, still
#    10..25 since step 1 did not add/remove paragraphs) with the new
#    6-paragraph block.
# ---------------------------------------------------------------------------
$startPara = $d.Paragraphs.Item(10)
$endPara = $d.Paragraphs.Item(25)
$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$blockBody = '<w:body>' `
  + '<w:p><w:r><w:lastRenderedPageBreak/><w:t>{</w:t></w:r><w:r><w:t>#</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>debut_offre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' `
  + '<w:p><w:pPr><w:pStyle w:val="Titre1"/></w:pPr><w:r><w:t>Prix: {Prix}</w:t></w:r></w:p>' `
  + '<w:p><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Titre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> {</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>titre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p>' `
  + '<w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fin_offre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p>' `
  + '<w:p><w:r><w:t>{</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>nom</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>}</w:t></w:r></w:p>' `
  + '<w:p><w:r><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>prenom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}</w:t></w:r></w:p>' `
  + '</w:body>'
$blockRange.InsertXML((New-FlatOpc $blockBody))

# ---------------------------------------------------------------------------
# 3. Footer: {nom} / {prenom} / {telephone} each wrapped with proofErr.
# ---------------------------------------------------------------------------
$ftr = $d.Sections.Item(1).Footers.Item(1)
$ftrPara = $ftr.Range.Paragraphs.Item(1)
$ftrBody = '<w:body><w:p><w:pPr><w:pStyle w:val="Pieddepage"/></w:pPr>' `
  + '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>{</w:t></w:r>' `
  + '<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>nom</w:t></w:r><w:proofErr w:type="gramEnd"/>' `
  + '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>}</w:t></w:r>' `
  + '<w:r><w:ptab w:relativeTo="margin" w:alignment="center" w:leader="none"/></w:r>' `
  + '<w:r><w:t>{</w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>prenom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' `
  + '<w:r><w:t>}</w:t></w:r>' `
  + '<w:r><w:ptab w:relativeTo="margin" w:alignment="right" w:leader="none"/></w:r>' `
  + '<w:r><w:t>{</w:t></w:r>' `
  + '<w:proofErr w:type="gramStart"/><w:r><w:t>telephone</w:t></w:r><w:proofErr w:type="gramEnd"/>' `
  + '<w:r><w:t>}</w:t></w:r>' `
  + '</w:p></w:body>'
$ftrPara.Range.InsertXML((New-FlatOpc $ftrBody))

# ---------------------------------------------------------------------------
# 4. Header: "{nom} {prenom}" wrapped with proofErr.
# ---------------------------------------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdrPara = $hdr.Range.Paragraphs.Item(1)
$hdrBody = '<w:body><w:p><w:pPr><w:pStyle w:val="En-tte"/></w:pPr>' `
  + '<w:r><w:t>{</w:t></w:r>' `
  + '<w:proofErr w:type="gramStart"/><w:r><w:t>nom</w:t></w:r><w:proofErr w:type="gramEnd"/>' `
  + '<w:r><w:t>} {</w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>prenom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' `
  + '<w:r><w:t>}</w:t></w:r>' `
  + '</w:p></w:body>'
$hdrPara.Range.InsertXML((New-FlatOpc $hdrBody))

Write-Host "Done."
